$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.552.57'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '1.660.28'
$ws.Range('E3').Value = '  -3.23%  '
$ws.Range('E4').Value = '  +0.73%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '214.60'
$ws.Range('E5').Value = '  -1.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '0.514'
$ws.Range('E6').Value = '  -1.22%  '
$ws.Range('E7').Value = '  +0.74%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range('D8').Value = '23.28'
$ws.Range('E8').Value = '  -2.50%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range('D9').Value = '0.258'
$ws.Range('E9').Value = '  -2.33%  '
$ws.Range('E10').Value = '  -1.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range('D11').Value = '0.0876'
$ws.Range('E11').Value = '  -2.09%  '
$ws.Range('D12').Value = '1.894.19'
$ws.Range('E12').Value = '  -3.32%  '
$ws.Range('D13').Value = '1.658.90'
$ws.Range('E13').Value = '  -3.16%  '
$ws.Range('E14').Value = '  -2.86%  '
$ws.Range('E15').Value = '  -2.69%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range('D16').Value = '65.77'
$ws.Range('E16').Value = '  -3.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range('D17').Value = '245.91'
$ws.Range('E17').Value = '  +1.28%  '
$ws.Range('D18').Value = '27.540.69'
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('D19').Value = '0.0₃0729'
$ws.Range('E19').Value = '  -2.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range('D20').Value = '7.47'
$ws.Range('E20').Value = '  -6.12%  '
$ws.Range('E21').Value = '  +0.75%  '
$ws.Range('E22').Value = '  -2.95%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '9.32'
$ws.Range('E23').Value = '  -3.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range('D24').Value = '2.02'
$ws.Range('E24').Value = '  -4.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range('D25').Value = '145.93'
$ws.Range('E25').Value = '  -1.60%  '
$ws.Range('E26').Value = '  -5.05%  '
$ws.Range('E27').Value = '  -2.36%  '
$ws.Range('E28').Value = '  +0.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range('D29').Value = '0.111'
$ws.Range('E29').Value = '  -2.07%  '
$ws.Range('E30').Value = '  +5.09%  '
$ws.Range('E31').Value = '  -0.79%  '
$ws.Range('E32').Value = '  -2.50%  '
$ws.Range('D33').Value = '1.438.07'
$ws.Range('E33').Value = '  -7.27%  '
$ws.Range('E34').Value = '  -5.13%  '
$ws.Range('E35').Value = '  -8.05%  '
$ws.Range('E36').Value = '  -0.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range('D37').Value = '0.929'
$ws.Range('E37').Value = '  -3.59%  '
$ws.Range('E38').Value = '  -5.74%  '
$ws.Range('E39').Value = '  -2.66%  '
$ws.Range('E40').Value = '  -1.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range('D41').Value = '69.09'
$ws.Range('E41').Value = '  -2.75%  '
$ws.Range('E42').Value = '  +0.87%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '5.42'
$ws.Range('E43').Value = '  -7.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range('D44').Value = '0.791'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').Value = '1.802.30'
$ws.Range('E45').Value = '  -3.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '2.20'
$ws.Range('E46').Value = '  -3.61%  '
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '88.71'
$ws.Range('E48').Value = '  -2.87%  '
$ws.Range('E49').Value = '  +4.61%  '
$ws.Range('E50').Value = '  -4.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range('D51').Value = '7.82'
$ws.Range('E51').Value = '  -5.69%  '
